$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 20 - Content Placeholder 2 - last bullet:
#    "Easily indefinitely scaling out because of decentralizing"
#      -> "Easily, indefinitely scaling out because of decentralizing"
# ---------------------------------------------------------------------------
$s20 = $p.Slides.Item(20)
$shp20 = $s20.Shapes.Item(2)
$tr20 = $shp20.TextFrame.TextRange
$run20 = $tr20.Characters(185, 57)
$run20.Text = "Easily, indefinitely scaling out because of decentralizing"

# ---------------------------------------------------------------------------
# 2) Slide 24 - Content Placeholder 2 - 2nd bullet (lvl 1):
#    "By chained read/write." -> "By chaining read/write."
# ---------------------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$shp24 = $s24.Shapes.Item(2)
$tr24 = $shp24.TextFrame.TextRange
$run24 = $tr24.Characters(68, 22)
$run24.Text = "By chaining read/write."

# ---------------------------------------------------------------------------
# 3) Slide 26 - Content Placeholder 2 - 3rd bullet (lvl 1):
#    "Optimizing by memory buffer" -> "Optimize by memory buffer"
# ---------------------------------------------------------------------------
$s26 = $p.Slides.Item(26)
$shp26 = $s26.Shapes.Item(2)
$tr26 = $shp26.TextFrame.TextRange
$run26 = $tr26.Characters(98, 27)
$run26.Text = "Optimize by memory buffer"

# ---------------------------------------------------------------------------
# 4) Slide 8 - Content Placeholder 2 - bullet "Small number of node in tree...":
#    split the run so that "number " becomes "& unchanged number " (inserted
#    phrase keeps its own run), giving:
#      "Small " | "& unchanged number " | "of node in tree, can be cache in
#      memory for fast "
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shp8 = $s8.Shapes.Item(2)
$tr8 = $shp8.TextFrame.TextRange
$mid8 = $tr8.Characters(104, 7)
$mid8.Text = "& unchanged number "
